# Generate functional Comet-Farm upload XML
#
# The "continue from previous year" block of the scenario sheet (rows 48-57,
# years 2020-2029) was missing the crop name and had placeholder tillage /
# nitrogen-application values. Fill in the crop, switch the tillage type to
# "No Tillage", and drop the N application amount to 5 so the generated
# upload XML is actually functional.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario")

for ($row = 48; $row -le 57; $row++) {
    $ws.Range("B$row").Value = "Alfalfa"
    $ws.Range("I$row").Value = "No Tillage"
    $ws.Range("M$row").Value = 5
}

# Leave the sheet scrolled/selected where the edit finished, as Excel would.
$ws.Range("N48").Select()
